$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata_Mean_results")

# Starting layout of "testdata_Mean_results":
#   A area | B value_sum | C value_count | D stdev | E value |
#   F lowercl | G uppercl | H confidence | I statistic | J method
# Rows 2-4  -> the "95%"  confidence-interval block (Area1, Area2, No grouping)
# Rows 5-7  -> the "99.8%" confidence-interval block (Area1, Area2, No grouping) - duplicate areas

# --- capture the 99.8% CI bounds (rows 5-7) before we start rearranging rows/columns ---
$lo99_Area1 = 18.216705294788838
$hi99_Area1 = 89.191984705211183
$lo99_Area2 = 3859.0770997295967
$hi99_Area2 = 7498.8488158259588
$lo99_Tot   = 1801.9535385474737
$hi99_Tot   = 6094.274769144833

# --- drop the duplicate 99.8% rows; only one row per area remains ---
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()

# --- drop the "confidence" column (col H); it will be replaced by two dedicated CI columns ---
$ws.Columns.Item(8).Delete()
# Layout now: A area | B value_sum | C value_count | D stdev | E value |
#             F lowercl | G uppercl | H statistic | I method

# --- insert two fresh columns for the 99.8% CI bounds, right after the 95% CI bounds ---
$ws.Columns.Item(8).Insert()
$ws.Columns.Item(8).Insert()
# Layout now: A area | B value_sum | C value_count | D stdev | E value |
#             F lowercl | G uppercl | H (new) | I (new) | J statistic | K method

# --- populate the new columns with the 99.8% CI bounds captured above ---
$ws.Cells.Item(2,8).Value = $lo99_Area1
$ws.Cells.Item(2,9).Value = $hi99_Area1

$ws.Cells.Item(3,8).Value = $lo99_Area2
$ws.Cells.Item(3,9).Value = $hi99_Area2

$ws.Cells.Item(4,8).Value = $lo99_Tot
$ws.Cells.Item(4,9).Value = $hi99_Tot

# --- rename headers: the two new columns become the 99.8% bounds, the two ---
# --- pre-existing CI columns become the 95% bounds ---
$ws.Cells.Item(1,8).Value = "lower99_8cl"
$ws.Cells.Item(1,9).Value = "upper99_8cl"
$ws.Cells.Item(1,6).Value = "lower95_0cl"
$ws.Cells.Item(1,7).Value = "upper95_0cl"

# --- match the saved selection state ---
$ws.Range("J1:J1048576").Select()
